$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H10").Value = 935547
$ws.Range("I10").Value = 0.547273
$ws.Range("J10").Value = 1120.82
$ws.Range("K10").Value = 58471.5
$ws.Range("L10").Value = 58630.6
$ws.Range("M10").Value = 58294.1
$ws.Range("N10").Value = 58630.5
$ws.Range("O10").Value = 58607.9
$ws.Range("P10").Value = 58460.6
$ws.Range("Q10").Value = 533.979
$ws.Range("R10").Value = 676.494
$ws.Range("S10").Value = 114.793
$ws.Range("T10").Value = 676.487
$ws.Range("U10").Value = 565.553
$ws.Range("V10").Value = 543.141
$ws.Range("W10").Value = 28.3036
$ws.Range("X10").Value = 28.575
$ws.Range("Y10").Value = 28.2279
$ws.Range("Z10").Value = 28.5428
$ws.Range("AA10").Value = 28.3525
$ws.Range("AB10").Value = 28.3117
$ws.Range("AD10").Value = 100
$ws.Range("AE10").Value = 92.02569444444444
$ws.Range("AG10").Value = 97.31447927607095
$ws.Range("AH10").Value = 91.30805528699688
